$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows at position 10 (pushes old rows 10-12 to 12-14) ---
$ws.Rows("10:11").Insert()

# --- Update recalculated numeric values in existing rows (2, 3, 4, 6, 7, 8) ---
# Row 2
$ws.Range("C2").Value = -1.191779160026281
$ws.Range("D2").Value = 9747.457259750538
$ws.Range("E2").Value = 17866.11560270452
$ws.Range("F2").Value = -0.4551661381186431
$ws.Range("G2").Value = 20544.39727538676
$ws.Range("H2").Value = 10278.90184040628

# Row 3
$ws.Range("C3").Value = -0.7811135038951063
$ws.Range("D3").Value = 8855.464493685831
$ws.Range("E3").Value = 16105.62345992551
$ws.Range("F3").Value = -0.07832636867183637
$ws.Range("G3").Value = 17290.32346778754
$ws.Range("H3").Value = 8779.115747629901

# Row 4
$ws.Range("F4").Value = -0.38354761714424

# Row 6
$ws.Range("C6").Value = 0.160688677823154
$ws.Range("D6").Value = 5589.49880952381
$ws.Range("E6").Value = 11055.87917927197
$ws.Range("F6").Value = -0.01202178513521324
$ws.Range("G6").Value = 18678.59998461572
$ws.Range("H6").Value = 8265.572647368421

# Row 7
$ws.Range("C7").Value = 0.2667827799071745
$ws.Range("D7").Value = 5683.118419047614
$ws.Range("E7").Value = 10333.51492798416
$ws.Range("F7").Value = 0.37054827230222
$ws.Range("G7").Value = 15032.85774066459
$ws.Range("H7").Value = 6483.545771421052

# Row 8
$ws.Range("C8").Value = 0.4328762177266555
$ws.Range("D8").Value = 4559.395693542933
$ws.Range("E8").Value = 9088.049273469329
$ws.Range("F8").Value = 0.1824812057294804
$ws.Range("G8").Value = 16160.67192591208
$ws.Range("H8").Value = 6996.381219108363

# --- Update recalculated inline-string (array) values in existing rows ---
# Row 2
$ws.Range("I2").Value = "[-5.36936080e+08 -4.84330921e+07 -6.09596707e+07 -3.33792368e+08`n -1.13024009e+09]"
$ws.Range("J2").Value = "[-1.86960885 -0.19760905 -0.61686118  0.56700031 -0.15875193]"
$ws.Range("K2").Value = "[ -9896.10240574  -6241.48377845  -5617.34012575 -11961.72481999`n -17677.85807211]"

# Row 3
$ws.Range("I3").Value = "[-3.76445292e+08 -3.31225718e+07 -5.31547431e+07 -3.41162024e+08`n -6.90891797e+08]"
$ws.Range("J3").Value = "[-1.01187959  0.18097544 -0.40984752  0.5574403   0.29167952]"
$ws.Range("K3").Value = "[ -8565.76416386  -4887.84968104  -5681.0247185  -11750.18556298`n -13010.75461177]"

# Row 6
$ws.Range("I6").Value = "[-5.06789192e+08 -1.19083116e+07 -2.57607087e+07 -5.28089041e+08`n -6.71903234e+08]"
$ws.Range("J6").Value = "[-1.70849139  0.7055422   0.31673696  0.3149562   0.3111471 ]"
$ws.Range("K6").Value = "[ -9917.2135      -2332.8615      -3178.4635     -12616.44315789`n -13282.88157895]"

# Row 7
$ws.Range("I7").Value = "[-2.88786445e+08 -8.02524460e+06 -1.85820434e+07 -2.84849348e+08`n -5.29690978e+08]"
$ws.Range("J7").Value = "[-0.5433944   0.80155912  0.50713998  0.63048982  0.45694685]"
$ws.Range("K7").Value = "[ -7037.863755    -2096.66447     -2724.67799     -9098.60597895`n -11459.91666316]"

# Row 8
$ws.Range("I8").Value = "[-4.08541864e+08 -7.83030482e+06 -2.84367527e+07 -2.70481201e+08`n -5.90546463e+08]"
$ws.Range("J8").Value = "[-1.18341696  0.80637941  0.24575903  0.64912836  0.3945562 ]"
$ws.Range("K8").Value = "[ -7910.78035296  -2054.01274835  -3212.11590815  -8928.71489846`n -12876.28218763]"

# --- Populate new row 10 (LightGBM) ---
$ws.Range("A10").Value = "Dataset 1"
$ws.Range("B10").Value = "LightGBM"
$ws.Range("C10").Value = -0.5062336308909425
$ws.Range("D10").Value = 9402.940989458264
$ws.Range("E10").Value = 14810.77939627252
$ws.Range("F10").Value = -0.7411596952515552
$ws.Range("G10").Value = 19530.44383360109
$ws.Range("H10").Value = 10822.51912859489
$ws.Range("I10").Value = "[-2.76605469e+08 -5.80258760e+07 -1.58249892e+08 -5.34481653e+08`n -8.79828292e+08]"
$ws.Range("J10").Value = "[-0.47829421 -0.43481061 -3.19733412  0.30666363  0.09797685]"
$ws.Range("K10").Value = "[-10018.5647826   -5643.43607803  -9006.64071117 -11911.40488544`n -17532.54918574]"

# --- Populate new row 11 (CatBoost) ---
$ws.Range("A11").Value = "Dataset 1"
$ws.Range("B11").Value = "CatBoost"
$ws.Range("C11").Value = 0.4032481595475715
$ws.Range("D11").Value = 4812.576743866931
$ws.Range("E11").Value = 9322.419174141065
$ws.Range("F11").Value = 0.3522514500934445
$ws.Range("G11").Value = 14858.62242451191
$ws.Range("H11").Value = 6530.750213051405
$ws.Range("I11").Value = "[-2.70000843e+08 -1.13316167e+07 -2.45974495e+07 -1.65618951e+08`n -6.32344441e+08]"
$ws.Range("J11").Value = "[-0.44299636  0.71980218  0.34759062  0.785157    0.35170381]"
$ws.Range("K11").Value = "[ -5841.93850766  -2563.5335398   -3167.97967961  -8365.2116056`n -12715.08773258]"

